$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "En esta sección..." paragraph: merge " para la verificación medi" +
#    "ante las pruebas del " into a single run, and drop the _GoBack bookmark
#    that used to sit right after it (it gets relocated further down below).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("verificación medi" + "ante las pruebas del ", $true, $false, $false, $false, $false, $true, 1, $false, "verificación mediante las pruebas del ", 2) | Out-Null

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2) "...muestre 0 RPM." -> "...muestre 0 +/- 10 RPM."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("en la pantalla de LCD  muestre 0 RPM.", $true, $false, $false, $false, $false, $true, 1, $false, "en la pantalla de LCD  muestre 0 +/- 10 RPM.", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) "...muestre 3000 RPM." -> "...muestre 3000  +/- 10 RPM."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("muestre 3000 RPM.", $true, $false, $false, $false, $false, $true, 1, $false, "muestre 3000  +/- 10 RPM.", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4) "...muestre 1500 +/- 50 RPM y déjelo con dicho valor." ->
#    "...muestre 1500 +/- 100 RPM y déjelo con dicho valor."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("muestre 1500 +/- 50 RPM", $true, $false, $false, $false, $false, $true, 1, $false, "muestre 1500 +/- 100 RPM", 2) | Out-Null

# ---------------------------------------------------------------------------
# 5) Fill in the trailing empty bulleted paragraph and append five more new
#    bulleted paragraphs describing the remaining test steps.
# ---------------------------------------------------------------------------

# Locate the last (empty) "Prrafodelista" bullet paragraph.
$lastIndex = $d.Paragraphs.Count
$emptyPara = $d.Paragraphs.Item($lastIndex - 2)
Write-Host "Target empty paragraph text: [" $emptyPara.Range.Text "]"

$emptyPara.Range.Text = "Mueva el potenciómetro hasta un valor de 2500 +/- 100, y verifique que el valor de RPM iguale al setpoint con una diferencia de +/- 100"

$cur = $d.Paragraphs.Item($lastIndex - 2)
$r = $cur.Range
$r2 = $r.Duplicate
$r2.Collapse(0)
$r2.InsertParagraphAfter()

$cur = $d.Paragraphs.Item($lastIndex - 1)
$cur.Range.Text = "Mueva el potenciómetro hasta un valor de 500 +/- 100, y verifique que el valor de RPM iguale al setpoint con una diferencia de +/- 100"

$r = $cur.Range
$r2 = $r.Duplicate
$r2.Collapse(0)
$r2.InsertParagraphAfter()

$cur = $d.Paragraphs.Item($lastIndex)
$cur.Range.Text = "Aplique una perturbación a la velocidad del motor aplicando alguna fuerza de oposición o mayor carga al motor y verifique que después de la perturbación la velocidad iguale al setpoint con una diferencia de +/- 100."

$r = $cur.Range
$r2 = $r.Duplicate
$r2.Collapse(0)
$r2.InsertParagraphAfter()

$cur = $d.Paragraphs.Item($lastIndex + 1)
$cur.Range.Text = "Presione cualquiera de los dos switches y verifique que el motor y el LED1 verde se apaguen, y verifique también en el LCD que los valores de trabajo de ciclo y las RPM estén dando una valor de cero."

$r = $cur.Range
$r2 = $r.Duplicate
$r2.Collapse(0)
$r2.InsertParagraphAfter()

$cur = $d.Paragraphs.Item($lastIndex + 2)
$cur.Range.Text = "Vuelva a presionar los dos botones y verifique que el valor de RPM iguale al setpoint con una diferencia de +/- 100"

$r = $cur.Range
$r2 = $r.Duplicate
$r2.Collapse(0)
$r2.InsertParagraphAfter()

$cur = $d.Paragraphs.Item($lastIndex + 3)
$cur.Range.Text = "Apague el motor presionando cualquiera de los dos switches de la tarjeta."

# ---------------------------------------------------------------------------
# 6) Re-insert the _GoBack bookmark inside the "Presione cualquiera..."
#    paragraph, right before "en el LCD que los valores...".
# ---------------------------------------------------------------------------
$matchRange = $d.Content
$found = $matchRange.Find.Execute("en el LCD que los valores de trabajo", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Host "Bookmark anchor found:" $found " at " $matchRange.Start "-" $matchRange.End
if ($found) {
    $bmRange = $d.Range($matchRange.Start, $matchRange.Start)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}

Write-Host "Done"
